# "minor text and link updates"
#
# 1) Slide 3 ("Advantages"): tweak the last bullet under "Liability:" to
#    call out "repetitive" tasks specifically.
# 2) The deck's cached "datetimeFigureOut" field text in every slide layout
#    rolled from 2/12/2026 to 2/13/2026 (PowerPoint recomputes these cached
#    captions whenever the file is saved on a later day).

$p = $ppt.ActivePresentation

# --- 1. Fix the bullet text on the "Advantages" slide -----------------
$advantagesSlide = $p.Slides.Item(3)
$contentShape = $advantagesSlide.Shapes.Item(2)
$tr = $contentShape.TextFrame.TextRange
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text -eq "Less productive than an agent for fixed tasks") {
        $para.Runs(1).Text = "Less productive than an agent for fixed repetitive tasks"
        break
    }
}

# --- 2. Refresh the cached date placeholder text on every slide layout -
$oldDateText = "2/12/2026"
$newDateText = "2/13/2026"

$master = $p.Designs.Item(1).SlideMaster
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    foreach ($shape in $layout.Shapes) {
        $isDatePlaceholder = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $shape.HasTextFrame) {
            $dateRange = $shape.TextFrame.TextRange
            if ($dateRange.Text -eq $oldDateText) {
                $dateRange.Text = $newDateText
            }
        }
    }
}
